# Auto-generated edit script: applies numeric corrections to the
# "currentAveragePrice*"/"LevePrice*"/"LeveProfit*" columns (H-N) of several
# leve-profit worksheets, per the scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2100
$ws.Range("I70").Value = 1375
$ws.Range("K70").Value = 4125
$ws.Range("M70").Value = -3855
$ws.Range("H73").Value = 2100
$ws.Range("I73").Value = 1375
$ws.Range("K73").Value = 4125
$ws.Range("M73").Value = -3189
$ws.Range("H113").Value = 40015824
$ws.Range("I113").Value = 50003970
$ws.Range("J113").Value = 63241
$ws.Range("K113").Value = 50003970
$ws.Range("L113").Value = 63241
$ws.Range("M113").Value = -50000716
$ws.Range("N113").Value = -69749
$ws.Range("H116").Value = 3424.875
$ws.Range("J116").Value = 3649.75
$ws.Range("L116").Value = 3649.75
$ws.Range("N116").Value = -10533.75
$ws.Range("H132").Value = 1401.9762
$ws.Range("I132").Value = 1249.8108
$ws.Range("K132").Value = 3749.4324
$ws.Range("M132").Value = -1219.4324
$ws.Range("H134").Value = 155956
$ws.Range("J134").Value = 155956
$ws.Range("L134").Value = 155956
$ws.Range("N134").Value = -166096
$ws.Range("H138").Value = 3049.6177
$ws.Range("I138").Value = 2042.25
$ws.Range("J138").Value = 3359.577
$ws.Range("K138").Value = 6126.75
$ws.Range("L138").Value = 10078.731
$ws.Range("M138").Value = -986.75
$ws.Range("N138").Value = -20358.731

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6067.311
$ws.Range("I32").Value = 6067.311
$ws.Range("K32").Value = 6067.311
$ws.Range("M32").Value = -5780.311
$ws.Range("H45").Value = 138830.6
$ws.Range("I45").Value = 186132.64
$ws.Range("K45").Value = 186132.64
$ws.Range("M45").Value = -185755.64
$ws.Range("H74").Value = 3517.3704
$ws.Range("I74").Value = 2598.76
$ws.Range("K74").Value = 2598.76
$ws.Range("M74").Value = -1724.76
$ws.Range("H77").Value = 3517.3704
$ws.Range("I77").Value = 2598.76
$ws.Range("K77").Value = 12993.8
$ws.Range("M77").Value = -8625.800000000001
$ws.Range("H110").Value = 6196
$ws.Range("I110").Value = 3861.5557
$ws.Range("K110").Value = 3861.5557
$ws.Range("M110").Value = -1816.5557
$ws.Range("H125").Value = 137500
$ws.Range("J125").Value = 137500
$ws.Range("L125").Value = 137500
$ws.Range("N125").Value = -147340
$ws.Range("H135").Value = 64516.5
$ws.Range("J135").Value = 64516.5
$ws.Range("L135").Value = 64516.5
$ws.Range("N135").Value = -74656.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2283.2727
$ws.Range("I94").Value = 1889.5555
$ws.Range("K94").Value = 1889.5555
$ws.Range("M94").Value = -1438.5555
$ws.Range("H105").Value = 1943.8572
$ws.Range("I105").Value = 1804.7646
$ws.Range("K105").Value = 1804.7646
$ws.Range("M105").Value = -57.76459999999997
$ws.Range("H134").Value = 2259.07
$ws.Range("I134").Value = 1974.9434
$ws.Range("K134").Value = 5924.8302
$ws.Range("M134").Value = -3389.8302
$ws.Range("H140").Value = 78994
$ws.Range("J140").Value = 78994
$ws.Range("L140").Value = 78994
$ws.Range("N140").Value = -89354

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H58").Value = 9152.166999999999
$ws.Range("I58").Value = 7194
$ws.Range("J58").Value = 11599.875
$ws.Range("K58").Value = 7194
$ws.Range("L58").Value = 11599.875
$ws.Range("M58").Value = -6991
$ws.Range("N58").Value = -12005.875
$ws.Range("H64").Value = 46837.625
$ws.Range("J64").Value = 49922.145
$ws.Range("L64").Value = 49922.145
$ws.Range("N64").Value = -50418.145
$ws.Range("H67").Value = 46837.625
$ws.Range("J67").Value = 49922.145
$ws.Range("L67").Value = 49922.145
$ws.Range("N67").Value = -51638.145
$ws.Range("H99").Value = 4958.4707
$ws.Range("I99").Value = 4363.364
$ws.Range("K99").Value = 4363.364
$ws.Range("M99").Value = -2865.364
$ws.Range("H126").Value = 4958.4707
$ws.Range("I126").Value = 4363.364
$ws.Range("K126").Value = 13090.092
$ws.Range("M126").Value = -10620.092
$ws.Range("H132").Value = 4450.579
$ws.Range("I132").Value = 2860.923
$ws.Range("K132").Value = 8582.769
$ws.Range("M132").Value = -6052.769
$ws.Range("H136").Value = 9152.166999999999
$ws.Range("I136").Value = 7194
$ws.Range("J136").Value = 11599.875
$ws.Range("K136").Value = 21582
$ws.Range("L136").Value = 34799.625
$ws.Range("M136").Value = -19032
$ws.Range("N136").Value = -39899.625
$ws.Range("H141").Value = 115982.695
$ws.Range("J141").Value = 125980.09
$ws.Range("L141").Value = 125980.09
$ws.Range("N141").Value = -136340.09

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 999.625
$ws.Range("I75").Value = 1062.4
$ws.Range("K75").Value = 3187.2
$ws.Range("M75").Value = -2189.2
$ws.Range("H78").Value = 999.625
$ws.Range("I78").Value = 1062.4
$ws.Range("K78").Value = 9561.6
$ws.Range("M78").Value = -4569.6
$ws.Range("H103").Value = 532.5
$ws.Range("I103").Value = 631.6667
$ws.Range("J103").Value = 433.33334
$ws.Range("K103").Value = 1895.0001
$ws.Range("L103").Value = 1300.00002
$ws.Range("M103").Value = -1016.0001
$ws.Range("N103").Value = -3058.00002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 36188.8
$ws.Range("I62").Value = 36188.8
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 36188.8
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -35502.8
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 36188.8
$ws.Range("I65").Value = 36188.8
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 108566.4
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -105134.4
$ws.Range("H107").Value = 715.1
$ws.Range("I107").Value = 772
$ws.Range("K107").Value = 772
$ws.Range("M107").Value = 1148
$ws.Range("H122").Value = 2932.6667
$ws.Range("I122").Value = 2932.6667
$ws.Range("K122").Value = 8798.000100000001
$ws.Range("M122").Value = -6348.000100000001
$ws.Range("H132").Value = 3428.7222
$ws.Range("I132").Value = 2440.2727
$ws.Range("J132").Value = 4982
$ws.Range("K132").Value = 7320.8181
$ws.Range("L132").Value = 14946
$ws.Range("M132").Value = -4790.8181
$ws.Range("N132").Value = -20006

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 517500
$ws.Range("I38").Value = 1000000
$ws.Range("J38").Value = 35000
$ws.Range("K38").Value = 1000000
$ws.Range("L38").Value = 35000
$ws.Range("M38").Value = -999590
$ws.Range("N38").Value = -35820
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H93").Value = 17716.076
$ws.Range("I93").Value = 2581.9
$ws.Range("K93").Value = 2581.9
$ws.Range("M93").Value = -1333.9
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H132").Value = 13346
$ws.Range("I132").Value = 14940.363
$ws.Range("J132").Value = 7500
$ws.Range("K132").Value = 44821.089
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -42291.089
$ws.Range("N132").Value = -27560
$ws.Range("H137").Value = 84700
$ws.Range("J137").Value = 84700
$ws.Range("L137").Value = 84700
$ws.Range("N137").Value = -94900
$ws.Range("H140").Value = 93806.664
$ws.Range("J140").Value = 93806.664
$ws.Range("L140").Value = 93806.664
$ws.Range("N140").Value = -104166.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 63838.5
$ws.Range("J46").Value = 63838.5
$ws.Range("L46").Value = 63838.5
$ws.Range("N46").Value = -64300.5
$ws.Range("H51").Value = 32000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 32000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 32000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -33020
$ws.Range("H81").Value = 2488.125
$ws.Range("I81").Value = 1701.3334
$ws.Range("K81").Value = 3402.6668
$ws.Range("M81").Value = -2341.6668
$ws.Range("H84").Value = 2488.125
$ws.Range("I84").Value = 1701.3334
$ws.Range("K84").Value = 17013.334
$ws.Range("M84").Value = -11709.334
$ws.Range("H100").Value = 887.7143
$ws.Range("I100").Value = 713.25
$ws.Range("K100").Value = 1426.5
$ws.Range("M100").Value = -885.5
$ws.Range("H107").Value = 1220.0714
$ws.Range("J107").Value = 1629.8572
$ws.Range("L107").Value = 4889.571599999999
$ws.Range("N107").Value = -8729.571599999999
$ws.Range("H122").Value = 3319.7437
$ws.Range("I122").Value = 2320.3572
$ws.Range("K122").Value = 6961.071599999999
$ws.Range("M122").Value = -4511.071599999999
$ws.Range("H132").Value = 4758.0586
$ws.Range("I132").Value = 4147.3706
$ws.Range("K132").Value = 12442.1118
$ws.Range("M132").Value = -9912.111800000001
$ws.Range("H134").Value = 63838.5
$ws.Range("J134").Value = 63838.5
$ws.Range("L134").Value = 191515.5
$ws.Range("N134").Value = -196585.5
